$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete the old "Run 50" column (AZ, containing the 51st run of
# data). This shifts the old "Mean" column (BA) left into AZ, and Excel
# auto-updates the sheet dimension / row spans for us.
$ws.Range("AZ1").EntireColumn.Delete()

# Step 2: Rename header A1 ("Gen" -> "MaxFES")
$ws.Range("A1").Value = "MaxFES"

# Step 3: Update column A (Gen -> MaxFES values), fill in previously-blank
# cells in columns L, M, T, AA, AI, AO, AW for rows 3-14, and write the
# recomputed per-row Mean values into column AZ.
$ws.Range("AZ2").Value = 16849669363.57773
$ws.Range("A3").Value = 0.001
$ws.Range("AA3").Value = 13881840481.86671
$ws.Range("AI3").Value = 9592025027.336111
$ws.Range("AO3").Value = 12598924151.49422
$ws.Range("AW3").Value = 15798494631.99224
$ws.Range("AZ3").Value = 11928288544.96775
$ws.Range("L3").Value = 23962772539.81879
$ws.Range("M3").Value = 9269472906.527388
$ws.Range("T3").Value = 11695016925.11665
$ws.Range("A4").Value = 0.01
$ws.Range("AA4").Value = 13881840481.86671
$ws.Range("AI4").Value = 7365880370.89247
$ws.Range("AO4").Value = 6411278973.132987
$ws.Range("AW4").Value = 6826037254.240469
$ws.Range("AZ4").Value = 1408888558.096272
$ws.Range("L4").Value = 7488452042.367895
$ws.Range("M4").Value = 5570808799.077648
$ws.Range("T4").Value = 8593934958.700897
$ws.Range("A5").Value = 0.1
$ws.Range("AA5").Value = 16207016.4152271
$ws.Range("AI5").Value = 7080006.02948127
$ws.Range("AO5").Value = 27184247.64106239
$ws.Range("AW5").Value = 6567252.42206041
$ws.Range("AZ5").Value = 20704954.23218983
$ws.Range("L5").Value = 9122373.346950511
$ws.Range("M5").Value = 47187783.65968932
$ws.Range("T5").Value = 13998422.50891949
$ws.Range("A6").Value = 0.2
$ws.Range("AA6").Value = 31438.44719703
$ws.Range("AI6").Value = 13329.52547113
$ws.Range("AO6").Value = 43649.35354769
$ws.Range("AW6").Value = 19452.97099146
$ws.Range("AZ6").Value = 18163864.92064066
$ws.Range("L6").Value = 22765.44404874
$ws.Range("M6").Value = 55906.16961238
$ws.Range("T6").Value = 106095.41894152
$ws.Range("A7").Value = 0.3
$ws.Range("AA7").Value = 3574.89748924
$ws.Range("AI7").Value = 1601.75750491
$ws.Range("AO7").Value = 2267.47214651
$ws.Range("AW7").Value = 9841.8605252
$ws.Range("AZ7").Value = 18158783.60555524
$ws.Range("L7").Value = 10545.24014355
$ws.Range("M7").Value = 1186.29533644
$ws.Range("T7").Value = 9554.052393149999
$ws.Range("A8").Value = 0.4
$ws.Range("AA8").Value = 133.46526792
$ws.Range("AI8").Value = 18.07039871
$ws.Range("AO8").Value = 42.88266747
$ws.Range("AW8").Value = 834.44349991
$ws.Range("AZ8").Value = 18158056.12031406
$ws.Range("L8").Value = 119.9457822
$ws.Range("M8").Value = 826.6026701
$ws.Range("T8").Value = 221.90319375
$ws.Range("A9").Value = 0.5
$ws.Range("AA9").Value = 5.57705553
$ws.Range("AI9").Value = 0.27093676
$ws.Range("AO9").Value = 0.49901676
$ws.Range("AW9").Value = 6.25128782
$ws.Range("AZ9").Value = 18158012.78122651
$ws.Range("L9").Value = 0.86147911
$ws.Range("M9").Value = 12.56462333
$ws.Range("T9").Value = 4.33470343
$ws.Range("A10").Value = 0.6
$ws.Range("AA10").Value = 0.02740003
$ws.Range("AI10").Value = 0.00307543
$ws.Range("AO10").Value = 0.00081781
$ws.Range("AW10").Value = 0.18879636
$ws.Range("AZ10").Value = 18158012.18895491
$ws.Range("L10").Value = 0.09644709
$ws.Range("M10").Value = 0.38477654
$ws.Range("T10").Value = 0.0442095
$ws.Range("A11").Value = 0.7
$ws.Range("AA11").Value = 0.00110263
$ws.Range("AI11").Value = 0.00021075
$ws.Range("AO11").Value = 0.00004075
$ws.Range("AW11").Value = 0.0079665
$ws.Range("AZ11").Value = 18158012.17432236
$ws.Range("L11").Value = 0.00071043
$ws.Range("M11").Value = 0.00305992
$ws.Range("T11").Value = 0.00080375
$ws.Range("A12").Value = 0.8
$ws.Range("AA12").Value = 0.00001238
$ws.Range("AI12").Value = 0.00000639
$ws.Range("AO12").Value = 0.00003661
$ws.Range("AW12").Value = 0.00004629
$ws.Range("AZ12").Value = 18158012.1740493
$ws.Range("L12").Value = 0.0000875
$ws.Range("M12").Value = 0.0000407
$ws.Range("T12").Value = 0.00001207
$ws.Range("A13").Value = 0.9
$ws.Range("AA13").Value = 0.00000036
$ws.Range("AI13").Value = 0.00000018
$ws.Range("AO13").Value = 0.00000092
$ws.Range("AW13").Value = 0.0000007
$ws.Range("AZ13").Value = 18158012.17404454
$ws.Range("L13").Value = 0.00000062
$ws.Range("M13").Value = 0.00000038
$ws.Range("T13").Value = 0.00000079
$ws.Range("A14").Value = 1
$ws.Range("AA14").Value = 0.00000001
$ws.Range("AI14").Value = 0.00000001
$ws.Range("AO14").Value = 0.00000001
$ws.Range("AZ14").Value = 18158012.17404446
$ws.Range("L14").Value = 0.00000001
$ws.Range("M14").Value = 0.00000001
$ws.Range("T14").Value = 0.00000001
